$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ', '
    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }
    if ($systemParts.Count -gt 0) {
        $newParts = $systemParts + $otherParts
        $newText = $newParts -join ', '
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
Write-Host "Done"
